$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 97.5
$ws.Range("I2").Value = 98.333336
$ws.Range("K2").Value = 98.333336
$ws.Range("M2").Value = 14.666664

$ws.Range("H15").Value = 327.85715
$ws.Range("I15").Value = 327.85715
$ws.Range("K15").Value = 983.5714499999999
$ws.Range("M15").Value = -814.5714499999999

$ws.Range("H53").Value = 515.8333
$ws.Range("I53").Value = 549
$ws.Range("K53").Value = 549
$ws.Range("M53").Value = 88

$ws.Range("H69").Value = 5000
$ws.Range("J69").Value = 5000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16748

$ws.Range("H70").Value = 2000
$ws.Range("J70").Value = 2750
$ws.Range("L70").Value = 8250
$ws.Range("N70").Value = -8790

$ws.Range("H72").Value = 5000
$ws.Range("J72").Value = 5000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53736

$ws.Range("H73").Value = 2000
$ws.Range("J73").Value = 2750
$ws.Range("L73").Value = 8250
$ws.Range("N73").Value = -10122

$ws.Range("H80").Value = 800
$ws.Range("I80").Value = 800
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2400
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1402
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 800
$ws.Range("I83").Value = 800
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 7200
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -2208
$ws.Range("N83").ClearContents()

$ws.Range("H86").Value = 9499.5
$ws.Range("I86").Value = 8999
$ws.Range("K86").Value = 8999
$ws.Range("M86").Value = -7876

$ws.Range("H89").Value = 9499.5
$ws.Range("I89").Value = 8999
$ws.Range("K89").Value = 44995
$ws.Range("M89").Value = -39379

$ws.Range("H92").Value = 574.6667
$ws.Range("I92").Value = 362.25
$ws.Range("J92").Value = 999.5
$ws.Range("K92").Value = 362.25
$ws.Range("L92").Value = 999.5
$ws.Range("M92").Value = 885.75
$ws.Range("N92").Value = -3495.5

$ws.Range("H100").Value = 2738.8
$ws.Range("I100").Value = 1423.75
$ws.Range("K100").Value = 1423.75
$ws.Range("M100").Value = -882.75

$ws.Range("H111").Value = 4710.375
$ws.Range("I111").Value = 497.5
$ws.Range("J111").Value = 6114.6665
$ws.Range("K111").Value = 1492.5
$ws.Range("L111").Value = 18343.9995
$ws.Range("M111").Value = 1574.5
$ws.Range("N111").Value = -24477.9995

$ws.Range("H135").Value = 1214.7646
$ws.Range("J135").Value = 2092.8
$ws.Range("L135").Value = 18835.2
$ws.Range("N135").Value = -23905.2

$ws.Range("H138").Value = 6269.8
$ws.Range("I138").Value = 7783
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 23349
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = -18209
$ws.Range("N138").Value = -22280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2490.5806
$ws.Range("I45").Value = 2265.3462
$ws.Range("J45").Value = 3661.8
$ws.Range("K45").Value = 2265.3462
$ws.Range("L45").Value = 3661.8
$ws.Range("M45").Value = -1888.3462
$ws.Range("N45").Value = -4415.8

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7641.25
$ws.Range("I134").Value = 2956.8572
$ws.Range("K134").Value = 8870.571599999999
$ws.Range("M134").Value = -6335.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 293.05554
$ws.Range("I7").Value = 340.41666
$ws.Range("K7").Value = 340.41666
$ws.Range("M7").Value = -227.41666

$ws.Range("H122").Value = 1578.2727
$ws.Range("I122").Value = 1557.1
$ws.Range("K122").Value = 4671.299999999999
$ws.Range("M122").Value = -2221.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 2113.111
$ws.Range("I23").Value = 2003
$ws.Range("J23").Value = 2333.3333
$ws.Range("K23").Value = 6009
$ws.Range("L23").Value = 6999.999899999999
$ws.Range("M23").Value = -5774
$ws.Range("N23").Value = -7469.999899999999

$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 6000
$ws.Range("M100").Value = -5189

$ws.Range("H128").Value = 100000
$ws.Range("I128").Value = 100000
$ws.Range("K128").Value = 300000
$ws.Range("M128").Value = -295020

$ws.Range("H140").Value = 2092
$ws.Range("I140").Value = 2092
$ws.Range("K140").Value = 6276
$ws.Range("M140").Value = -1096

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 1516.5
$ws.Range("I22").Value = 1619.8
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1619.8
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -1090.8
$ws.Range("N22").Value = -2058

$ws.Range("H97").Value = 1786.1428
$ws.Range("I97").Value = 1431.75
$ws.Range("J97").Value = 2258.6667
$ws.Range("K97").Value = 1431.75
$ws.Range("L97").Value = 2258.6667
$ws.Range("M97").Value = -935.75
$ws.Range("N97").Value = -3250.6667

$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1370
$ws.Range("N113").ClearContents()

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 819

$ws.Range("H27").Value = 819

$ws.Range("H42").Value = 35512.5
$ws.Range("I42").Value = 26025
$ws.Range("K42").Value = 26025
$ws.Range("M42").Value = -25462

$ws.Range("H43").Value = 11506
$ws.Range("I43").Value = 11506
$ws.Range("K43").Value = 11506
$ws.Range("M43").Value = -11313

$ws.Range("H46").Value = 880
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 1000
$ws.Range("N46").Value = -1376

$ws.Range("H49").Value = 35512.5
$ws.Range("I49").Value = 26025
$ws.Range("K49").Value = 26025
$ws.Range("M49").Value = -25878

$ws.Range("H82").Value = 1999.8334
$ws.Range("I82").Value = 1999.8334
$ws.Range("K82").Value = 1999.8334
$ws.Range("M82").Value = -1638.8334

$ws.Range("H85").Value = 1999.8334
$ws.Range("I85").Value = 1999.8334
$ws.Range("K85").Value = 1999.8334
$ws.Range("M85").Value = -751.8334

$ws.Range("H100").Value = 4480
$ws.Range("I100").Value = 4666.6665
$ws.Range("J100").Value = 4200
$ws.Range("K100").Value = 4666.6665
$ws.Range("L100").Value = 4200
$ws.Range("M100").Value = -4125.6665
$ws.Range("N100").Value = -5282

$ws.Range("H134").Value = 49863
$ws.Range("J134").Value = 49863
$ws.Range("L134").Value = 49863
$ws.Range("N134").Value = -60003

$ws.Range("H136").Value = 7254.364
$ws.Range("I136").Value = 3571.1428
$ws.Range("K136").Value = 10713.4284
$ws.Range("M136").Value = -8163.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 3026
$ws.Range("I37").Value = 3026
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 3026
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2823
$ws.Range("N37").ClearContents()

$ws.Range("H110").Value = 27999.5
$ws.Range("J110").Value = 27999.5
$ws.Range("L110").Value = 27999.5
$ws.Range("N110").Value = -36179.5

$ws.Range("H122").Value = 402319
$ws.Range("I122").Value = 1001750
$ws.Range("J122").Value = 2698.3333
$ws.Range("K122").Value = 3005250
$ws.Range("L122").Value = 8094.999899999999
$ws.Range("M122").Value = -3002800
$ws.Range("N122").Value = -12994.9999

$ws.Range("H126").Value = 1125.75
$ws.Range("I126").Value = 1301
$ws.Range("K126").Value = 3903
$ws.Range("M126").Value = -1433

$ws.Range("H132").Value = 6319.3184
$ws.Range("I132").Value = 3130.7856
$ws.Range("J132").Value = 11899.25
$ws.Range("K132").Value = 9392.356800000001
$ws.Range("L132").Value = 35697.75
$ws.Range("M132").Value = -6862.356800000001
$ws.Range("N132").Value = -40757.75

$ws.Range("H136").Value = 1675.5555
$ws.Range("I136").Value = 1510
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 4530
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1980
$ws.Range("N136").Value = -14100
